$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# The handback has completed and is in sync with en-US: update the shared
# "Status" text everywhere it is shown (Overview summary columns, and the
# Status column on each language detail sheet).
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# zh-cn: handback just regenerated - refresh the handback datetime and
# clear the stale "handback file is not latest" error detail.
$zhcn.Range("K2").Value = "2016-08-13 23:05:16"
$zhcn.Range("P2").Value = ""

# de-de: handback just regenerated - refresh the handback datetime and
# clear the stale "handback file is not latest" error detail.
$dede.Range("K2").Value = "2016-08-13 23:05:26"
$dede.Range("P2").Value = ""

# Let the changed content re-flow the affected column widths.
$overview.Columns("E:F").AutoFit() | Out-Null
$zhcn.Columns("C:C").AutoFit() | Out-Null
$zhcn.Columns("P:P").AutoFit() | Out-Null
$dede.Columns("C:C").AutoFit() | Out-Null
$dede.Columns("P:P").AutoFit() | Out-Null
